$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings (column D) stay as text, matching the
# original inline-string storage, instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.432.04'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '1.843.86'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '264.54'
$ws.Range("E5").Value = '  -3.45%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.5196'
$ws.Range("E7").Value = '  -1.63%  '
$ws.Range("D8").Value = '0.3268'
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("D9").Value = '0.06796'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("E10").Value = '  -5.69%  '
$ws.Range("D11").Value = '0.7775'
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").Value = '0.07744'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = '1.823.94'
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").Value = '88.03'
$ws.Range("E14").Value = '  -2.30%  '
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("D16").Value = '0.9996'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '13.92'
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '0.000007953'
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").Value = '26.453.81'
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("D21").Value = '2.076.46'
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").Value = '4.622'
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("D23").Value = '9.574'
$ws.Range("E23").Value = '  -3.97%  '
$ws.Range("D24").Value = '5.999'
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = '143.72'
$ws.Range("D26").Value = '2.179'
$ws.Range("E26").Value = '  -8.21%  '
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("D29").Value = '111.96'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").Value = '4.170'
$ws.Range("E30").Value = '  -3.70%  '
$ws.Range("D31").Value = '4.124'
$ws.Range("E31").Value = '  -4.42%  '
$ws.Range("D32").Value = '0.08713'
$ws.Range("D33").Value = '0.04827'
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.7210'
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.131'
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("D36").Value = '2.842'
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").Value = '3.106'
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.231'
$ws.Range("E38").Value = '  -4.22%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01779'
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("D40").Value = '0.4871'
$ws.Range("E40").Value = '  -4.28%  '
$ws.Range("D41").Value = '0.9149'
$ws.Range("E41").Value = '  -2.60%  '
$ws.Range("D42").Value = '111.00'
$ws.Range("E42").Value = '  -4.69%  '
$ws.Range("D43").Value = '6.063'
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = '7.725'
$ws.Range("E45").Value = '  -3.38%  '
$ws.Range("D46").Value = '0.4173'
$ws.Range("E46").Value = '  -5.41%  '
$ws.Range("D47").Value = '0.05922'
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").Value = '9.057'
$ws.Range("E48").Value = '  -3.04%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1240'
$ws.Range("E49").Value = '  -6.64%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '35.03'
$ws.Range("E50").Value = '  -2.96%  '
$ws.Range("D51").Value = '0.8855'
$ws.Range("E51").Value = '  +0.66%  '

# Restore the default "Normal" style on column D so no stray number-format
# style index is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"
